$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.939.38"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.087.41"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.75"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.85"
$ws.Range("E6").Value = "  -2.28%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.085.77"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.38"
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.619.54"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.51"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000162"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.085.71"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.091.82"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.90"
$ws.Range("E19").Value = "  -2.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.49"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "348.23"
$ws.Range("E22").Value = "  +2.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.84"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.13"
$ws.Range("E25").Value = "  +1.94%  "
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0882"
$ws.Range("E29").Value = "  -2.59%  "
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  +1.53%  "
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.92"
$ws.Range("E33").Value = "  -6.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.79"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.97"
$ws.Range("E35").Value = "  +8.54%  "
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.19"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.96"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("E43").Value = "  +1.72%  "
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.388.91"
$ws.Range("E45").Value = "  +5.30%  "
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.128.04"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.962"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("E51").Value = "  -1.60%  "
